# Update the symbol/price list as scraped on Mon Dec 12 21:33:19 UTC 2022
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value.
# Using NumberFormat "@" (text) first so Excel keeps the exact string
# (preserves trailing zeros / leading zeros) instead of coercing to a number.
$updates = [ordered]@{
    "D2"  = "275.98"
    "D3"  = "21.13"
    "D4"  = "6.250"
    "D5"  = "0.06233"
    "D6"  = "3.547"
    "D7"  = "1.546"
    "D8"  = "6.554"
    "D10" = "0.1652"
    "D11" = "0.08288"
    "D12" = "0.03493"
    "D13" = "0.03125"
    "D15" = "3.762"
    "D16" = "0.001627"
    "D17" = "0.04686"
    "D19" = "0.006229"
    "D20" = "0.001066"
    "D21" = "0.0001496"
    "D24" = "0.01396"
    "D25" = "0.3291"
    "D26" = "0.1248"
    "D28" = "0.0002729"
    "D40" = "0.04747"
    "D41" = "0.005287"
    "E41" = "40CEJICEJIBestin24h"
    "D42" = "0.007034"
    "E42" = "41KickTokenKICK"
    "D43" = "0.1121"
    "D44" = "0.01135"
    "D45" = "0.00006216"
    "D47" = "0.7211"
    "D49" = "0.00001895"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
